$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("F2").Value = 25.69000000000058
$ws.Range("H2").Value = [double]"8.792373606958392e-09"
$ws.Range("I2").Value = [double]"8.792373606958392e-09"
$ws.Range("L2").Value = 46.5728414631709
$ws.Range("M2").Value = '[32.23372486717266, 60.91195805916913]'
$ws.Range("N2").Value = [double]"4.886800630821142e-08"
$ws.Range("O2").Value = [double]"4.886800630821142e-08"
$ws.Range("P2").Value = 1.427710775505272
$ws.Range("Q2").Value = '[1.0629212381515014, 1.7925003128590422]'
$ws.Range("R2").Value = [double]"5.108320433322433e-10"
$ws.Range("S2").Value = [double]"5.108320433322433e-10"
$ws.Range("T2").Value = 57.04584444803056
$ws.Range("U2").Value = '[48.06788677043703, 66.0238021256241]'
$ws.Range("V2").Value = 0
$ws.Range("W2").Value = 0
$ws.Range("X2").Value = 19.85253253253298
$ws.Range("Y2").Value = 18.36102102102143
$ws.Range("Z2").Value = 21.34404404404452

# Row 3
$ws.Range("F3").Value = 25.69000000000058
$ws.Range("H3").Value = [double]"1.5748069515098e-10"
$ws.Range("I3").Value = [double]"1.5748069515098e-10"
$ws.Range("L3").Value = 50.13902331139472
$ws.Range("M3").Value = '[35.59492813811012, 64.68311848467933]'
$ws.Range("N3").Value = [double]"1.236382662916924e-08"
$ws.Range("O3").Value = [double]"1.236382662916924e-08"
$ws.Range("P3").Value = 1.226447582482503
$ws.Range("Q3").Value = '[0.9119738433844251, 1.54092132158058]'
$ws.Range("R3").Value = [double]"5.608367104059653e-10"
$ws.Range("S3").Value = [double]"5.608367104059653e-10"
$ws.Range("T3").Value = 54.37877368578868
$ws.Range("U3").Value = '[46.16600941751654, 62.59153795406082]'
$ws.Range("V3").Value = 0
$ws.Range("W3").Value = 0
$ws.Range("X3").Value = 20.6754354354359
$ws.Range("Y3").Value = 19.38964964965008
$ws.Range("Z3").Value = 21.96122122122171

# Row 4
$ws.Range("F4").Value = 25.69000000000058
$ws.Range("H4").Value = [double]"2.589625369697757e-08"
$ws.Range("I4").Value = [double]"2.589625369697757e-08"
$ws.Range("L4").Value = 48.42157825414616
$ws.Range("M4").Value = '[30.400905776784185, 66.44225073150814]'
$ws.Range("N4").Value = [double]"2.308495029135216e-06"
$ws.Range("O4").Value = [double]"2.308495029135216e-06"
$ws.Range("P4").Value = 0.6352369529781168
$ws.Range("Q4").Value = '[0.25786846606042424, 1.0126054398958093]'
$ws.Range("R4").Value = 0.001461943992448367
$ws.Range("S4").Value = 0.001461943992448367
$ws.Range("T4").Value = 48.6343790685756
$ws.Range("U4").Value = '[39.38821005961305, 57.88054807753816]'
$ws.Range("V4").Value = [double]"8.282263763703668e-14"
$ws.Range("W4").Value = [double]"8.282263763703668e-14"
$ws.Range("X4").Value = 23.09271271271323
$ws.Range("Y4").Value = 21.54976976977025
$ws.Range("Z4").Value = 24.6356556556562

# Row 5
$ws.Range("B5").Value = 0
$ws.Range("F5").Value = 25.69000000000058
$ws.Range("H5").Value = [double]"3.897484601722567e-08"
$ws.Range("I5").Value = [double]"3.897484601722567e-08"
$ws.Range("L5").Value = 46.0052080473273
$ws.Range("M5").Value = '[30.068249010476542, 61.94216708417806]'
$ws.Range("N5").Value = [double]"5.889916725454469e-07"
$ws.Range("O5").Value = [double]"5.889916725454469e-07"
$ws.Range("P5").Value = 0.2956053147521933
$ws.Range("Q5").Value = '[-0.09434212172942225, 0.6855527512338089]'
$ws.Range("R5").Value = 0.1338044896403581
$ws.Range("S5").Value = 0.1338044896403581
$ws.Range("T5").Value = 55.63375633683128
$ws.Range("U5").Value = '[46.79630712077991, 64.47120555288265]'
$ws.Range("V5").Value = [double]"2.220446049250313e-16"
$ws.Range("W5").Value = [double]"2.220446049250313e-16"
$ws.Range("X5").Value = 24.48136136136191
$ws.Range("Y5").Value = 22.8869869869875
$ws.Range("Z5").Value = 26.07573573573632

# Row 6
$ws.Range("F6").Value = 23.3100000000002
$ws.Range("H6").Value = [double]"3.133219239614959e-10"
$ws.Range("I6").Value = [double]"3.133219239614959e-10"
$ws.Range("L6").Value = 51.12785472923869
$ws.Range("M6").Value = '[36.00093948137469, 66.2547699771027]'
$ws.Range("N6").Value = [double]"1.967248008938327e-08"
$ws.Range("O6").Value = [double]"1.967248008938327e-08"
$ws.Range("P6").Value = -0.2641579408423853
$ws.Range("Q6").Value = '[-0.5912106295043857, 0.06289474781961513]'
$ws.Range("R6").Value = 0.1107679696613701
$ws.Range("S6").Value = 0.1107679696613701
$ws.Range("T6").Value = 55.76002196706307
$ws.Range("U6").Value = '[47.192825273158846, 64.32721866096729]'
$ws.Range("V6").Value = 0
$ws.Range("W6").Value = 0
$ws.Range("X6").Value = 0.9800000000000111
$ws.Range("Y6").Value = -0.2333333333333327
$ws.Range("Z6").Value = 2.193333333333355

# Row 7
$ws.Range("F7").Value = 23.3100000000002
$ws.Range("H7").Value = [double]"1.447477723237611e-07"
$ws.Range("I7").Value = [double]"1.447477723237611e-07"
$ws.Range("L7").Value = 42.53727188471691
$ws.Range("M7").Value = '[26.897141137798975, 58.17740263163484]'
$ws.Range("N7").Value = [double]"1.846819479789374e-06"
$ws.Range("O7").Value = [double]"1.846819479789374e-06"
$ws.Range("P7").Value = -0.2138421425866923
$ws.Range("Q7").Value = '[-0.6415264277600778, 0.2138421425866932]'
$ws.Range("R7").Value = 0.3192945136149061
$ws.Range("S7").Value = 0.3192945136149061
$ws.Range("T7").Value = 51.44468696250733
$ws.Range("U7").Value = '[42.35771751609649, 60.531656408918174]'
$ws.Range("V7").Value = [double]"7.327471962526033e-15"
$ws.Range("W7").Value = [double]"7.327471962526033e-15"
$ws.Range("X7").Value = 0.7933333333333401
$ws.Range("Y7").Value = -0.7933333333333414
$ws.Range("Z7").Value = 2.380000000000022

# Row 8
$ws.Range("F8").Value = 23.3100000000002
$ws.Range("H8").Value = [double]"5.651766388226065e-09"
$ws.Range("I8").Value = [double]"5.651766388226065e-09"
$ws.Range("L8").Value = 50.9408680132377
$ws.Range("M8").Value = '[35.47406477323125, 66.40767125324416]'
$ws.Range("N8").Value = [double]"3.567886763811146e-08"
$ws.Range("O8").Value = [double]"3.567886763811146e-08"
$ws.Range("P8").Value = 0.09434212172942225
$ws.Range("Q8").Value = '[-0.27044741562434815, 0.45913165908319264]'
$ws.Range("R8").Value = 0.6049981134181504
$ws.Range("S8").Value = 0.6049981134181504
$ws.Range("T8").Value = 51.98518339694243
$ws.Range("U8").Value = '[42.47247184537703, 61.497894948507835]'
$ws.Range("V8").Value = [double]"2.375877272697835e-14"
$ws.Range("W8").Value = [double]"2.375877272697835e-14"
$ws.Range("X8").Value = 22.96000000000021
$ws.Range("Y8").Value = 21.60666666666686
$ws.Range("Z8").Value = 24.31333333333355
